$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 618, shifting existing rows 618-651 down to 619-652
$ws.Rows.Item(618).Insert()

# Populate the new row 618 with data
$ws.Cells.Item(618, 1).Value = 6
$ws.Cells.Item(618, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(618, 3).Value = "Metropolitana"
$ws.Cells.Item(618, 4).Value = 44931
$ws.Cells.Item(618, 5).Value = 13
$ws.Cells.Item(618, 6).Value = 100112044
$ws.Cells.Item(618, 7).Value = "Perejil"
$ws.Cells.Item(618, 8).Value = "Sin especificar"
$ws.Cells.Item(618, 9).Value = "Primera"
$ws.Cells.Item(618, 10).Value = 310
$ws.Cells.Item(618, 11).Value = 10000
$ws.Cells.Item(618, 12).Value = 11000
$ws.Cells.Item(618, 13).Value = 10484
$ws.Cells.Item(618, 14).Value = '$/docena de atados'
$ws.Cells.Item(618, 15).Value = "Región Metropolitana"
$ws.Cells.Item(618, 16).Value = 3495
$ws.Cells.Item(618, 17).Value = 3
$ws.Cells.Item(618, 18).Value = "Hortaliza"
